# Add a new "Clusterware Version" column (column V / 22) to the Oracle DB
# template, to the right of the existing "Pgsql Migrability" column (U / 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve the existing "Pgsql Migrability" header's wrap-text formatting
# (it is re-asserted here so it survives the round trip unchanged).
$existingLastHeader = $ws.Cells.Item(1, 21)
$existingLastHeader.WrapText = $true

# --- New header cell: "Clusterware Version" in column V (22), row 1.
$newHeader = $ws.Cells.Item(1, 22)
$newHeader.Value = "Clusterware Version"

# Match the look of the other header cells: bold font, centered both ways.
$newHeader.Font.Bold = $true
$newHeader.HorizontalAlignment = -4108   # xlCenter
$newHeader.VerticalAlignment = -4108     # xlCenter
$newHeader.WrapText = $false

# --- Give the new column a sensible width, same family as the other columns.
$ws.Columns.Item(22).ColumnWidth = 17.29

# --- Move the selection to the newly added header cell.
[void]$ws.Range("V1").Select()

Write-Output "Added 'Clusterware Version' column (V1)."
